$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 31; this shifts existing rows 31-37 down to 32-38
$ws.Rows.Item(31).Insert()

# Populate the new row 31 with data (same shape as the former row 31, with
# the date/volume/price columns updated to the new week's values)
$ws.Range("A31").Value = 5
$ws.Range("B31").Value = "Macroferia Regional de Talca"
$ws.Range("C31").Value = "Maule"
$ws.Range("D31").Value = 44748
$ws.Range("E31").Value = 7
$ws.Range("F31").Value = 100112043
$ws.Range("G31").Value = "Pepino dulce"
$ws.Range("H31").Value = "Cultivar IV Región"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 300
$ws.Range("K31").Value = 15000
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = 15000
$ws.Range("N31").Value = "$/bandeja 18 kilos"
$ws.Range("O31").Value = "Provincia de Limarí"
$ws.Range("P31").Value = 833
$ws.Range("Q31").Value = 18
$ws.Range("R31").Value = "Hortaliza"
